$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6, pushing old rows 6-17 down to 7-18
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with data (mirrors row 5's pattern: otu, count, state, citation)
$ws.Range("A6").Value = "Aus aus"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = "Alaska"
$ws.Range("I6").Value = "Smith, 1920. Bears on the coast. Jr. Chilly Waters. 0:0 pp0-40."

# Match row height used by similar data rows (same as row 5)
$ws.Rows.Item(6).RowHeight = 105

# Update dimension/view to reflect the new row count and refreshed selection
$null = $ws.Range("H6").Select()
